# Applies the "Elimina EC anteriores y se agregan nuevos, se modifica base de datos" edit:
# - Updates the header "VALOR MORA" / counts
# - Replaces the 3-row worker detail table with a new 19-row table
# - Keeps the footer (signature) block, now pushed further down

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header block updates (row 11: total "VALOR MORA"; row 13: counts)
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 929400
$ws.Range("C13").Value = 12
$ws.Range("F13").Value = 9

# ---------------------------------------------------------------------------
# 2) Grow the worker-detail table.
#    Before: header row 15, data rows 16-18 (18 = last/bottom-border row),
#    then a blank gap (19-22) and the signature footer at 23-24.
#    After:  header row 15, data rows 16-34 (34 = last/bottom-border row),
#    the same blank gap shifted to 35-38 and the footer at 39-40.
#    We insert 16 plain rows right above the old "last" row (row 18) so that
#    row keeps its special bottom-border formatting as the new last row (34),
#    and the footer block is pushed down automatically along with it.
# ---------------------------------------------------------------------------
$ws.Range("18:33").Insert()

# Re-apply the "regular" data-row formatting (copied from row 17, the plain
# style used by every non-last data row) onto the freshly inserted rows.
$ws.Range("B17:J17").Copy()
$ws.Range("B18:J33").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Fill in the new worker-detail table contents (rows 16-34)
# ---------------------------------------------------------------------------
$rows = @(
  @(16,'CC','30776395','YARIMA INES ALCALA SIERRA','2507',56940,1423500),
  @(17,'CC','45592337','GISELA MARGARITA BAENA QUINTANA','2507',56940,1423500),
  @(18,'CC','30774843','BEATRIZ PUELLO GONZALEZ','2507',56940,1423500),
  @(19,'CC','9099307','ARMANDO ENRIQUE SERNA CORREA','2507',40000,1000000),
  @(20,'CC','9099307','ARMANDO ENRIQUE SERNA CORREA','2506',40000,1000000),
  @(21,'CC','9099307','ARMANDO ENRIQUE SERNA CORREA','2505',40000,1000000),
  @(22,'CC','9099307','ARMANDO ENRIQUE SERNA CORREA','2504',40000,1000000),
  @(23,'CC','9099307','ARMANDO ENRIQUE SERNA CORREA','2503',40000,1000000),
  @(24,'CC','9099307','ARMANDO ENRIQUE SERNA CORREA','2502',40000,1000000),
  @(25,'CC','1050952239','MARYORY CARDONA CASTRO','2507',56940,781242),
  @(26,'CC','9293298','DUPERLEY NAVARRO MOSCOTE','2507',56940,908526),
  @(27,'CC','22999741','EDAISY DIAZ PUELLO','2507',56940,877803),
  @(28,'CC','7937269','PEDRO LUIS GUERRA DONADO','2507',56940,877803),
  @(29,'CC','33255442','SIRLYS FLOREZ HERNANDEZ','2205',40000,1000000),
  @(30,'CC','33255442','SIRLYS FLOREZ HERNANDEZ','2204',40000,1000000),
  @(31,'CC','33255442','SIRLYS FLOREZ HERNANDEZ','2203',40000,1000000),
  @(32,'CC','30765312','GLORIA MARIA SARMIENTO FERNANDEZ','2507',56940,877803),
  @(33,'CC','45496886','LILIANA ZORALLA BENT MENDEZ','2507',56940,877803),
  @(34,'CC','45591359','LILIANA INES CASTRO SIERRA','2507',56940,1423500)
)

foreach ($r in $rows) {
  $rowNum = $r[0]
  $ws.Range("B$rowNum").Value = $r[1]
  $ws.Range("C$rowNum").Value = $r[2]
  $ws.Range("D$rowNum").Value = $r[3]
  $ws.Range("E$rowNum").Value = $r[4]
  $ws.Range("F$rowNum").Value = $r[5]
  $ws.Range("G$rowNum").Value = $r[6]
}
